$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Every data cell in this sheet is stored as text (Coin/Link/Price/
# Volume columns all come from t="inlineStr" cells in the source file,
# even values that look numeric, e.g. "300.55" or "43.140.29", and the
# "  -0.56%  " style percentage strings). Excel normally auto-coerces a
# plain Range.Value assignment like "300.55" into a real number, so we
# temporarily force the edited range to Text format first; ClearFormats()
# afterwards drops that now-unneeded formatting so the cells keep their
# original (default) style.
$editRange = $ws.Range("B2:E51")
$editRange.NumberFormat = "@"

$ws.Range('D2').Value = '43.140.29'
$ws.Range('E2').Value = '  -0.56%  '
$ws.Range('D3').Value = '2.305.39'
$ws.Range('E3').Value = '  -0.63%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').Value = '300.55'
$ws.Range('E5').Value = '  -0.69%  '
$ws.Range('D6').Value = '98.03'
$ws.Range('E6').Value = '  -3.56%  '
$ws.Range('E7').Value = '  +3.27%  '
$ws.Range('E8').Value = '  +0.01%  '
$ws.Range('D9').Value = '0.516'
$ws.Range('E9').Value = '  -0.58%  '
$ws.Range('D10').Value = '35.78'
$ws.Range('E10').Value = '  -2.69%  '
$ws.Range('E11').Value = '  -0.51%  '
$ws.Range('E12').Value = '  +0.33%  '
$ws.Range('D13').Value = '18.04'
$ws.Range('E13').Value = '  -3.98%  '
$ws.Range('E14').Value = '  -1.59%  '
$ws.Range('D15').Value = '2.663.58'
$ws.Range('E15').Value = '  -0.42%  '
$ws.Range('D16').Value = '2.301.53'
$ws.Range('E16').Value = '  -0.27%  '
$ws.Range('D17').Value = '0.788'
$ws.Range('E17').Value = '  -2.60%  '
$ws.Range('D18').Value = '43.028.96'
$ws.Range('E18').Value = '  -0.50%  '
$ws.Range('D19').Value = '13.32'
$ws.Range('E19').Value = '  +6.45%  '
$ws.Range('E20').Value = '  +0.63%  '
$ws.Range('E21').Value = '  -2.06%  '
$ws.Range('D22').Value = '68.54'
$ws.Range('E22').Value = '  +0.27%  '
$ws.Range('D23').Value = '238.20'
$ws.Range('E23').Value = '  +0.27%  '
$ws.Range('E24').Value = '  -3.36%  '
$ws.Range('E25').Value = '  -0.37%  '
$ws.Range('E26').Value = '  -1.97%  '
$ws.Range('D27').Value = '24.89'
$ws.Range('E27').Value = '  -1.90%  '
$ws.Range('D28').Value = '168.05'
$ws.Range('E28').Value = '  -0.68%  '
$ws.Range('E29').Value = '  -0.84%  '
$ws.Range('E30').Value = '  -7.65%  '
$ws.Range('D31').Value = '32.92'
$ws.Range('E31').Value = '  -6.96%  '
$ws.Range('D32').Value = '5.17'
$ws.Range('E32').Value = '  +1.91%  '
$ws.Range('E33').Value = '  +0.03%  '
$ws.Range('B34').Value = 'Celestia'
$ws.Range('C34').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D34').Value = '18.13'
$ws.Range('E34').Value = '  +0.89%  '
$ws.Range('B35').Value = 'RenderToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D35').Value = '4.76'
$ws.Range('E35').Value = '  +0.17%  '
$ws.Range('D36').Value = '2.41'
$ws.Range('E36').Value = '  -0.39%  '
$ws.Range('D37').Value = '0.0688'
$ws.Range('E37').Value = '  -1.80%  '
$ws.Range('E38').Value = '  +0.45%  '
$ws.Range('E39').Value = '  -1.32%  '
$ws.Range('E40').Value = '  +0.87%  '
$ws.Range('E41').Value = '  -3.98%  '
$ws.Range('D42').Value = '2.009.92'
$ws.Range('E42').Value = '  +0.49%  '
$ws.Range('E43').Value = '  -1.27%  '
$ws.Range('D44').Value = '2.17'
$ws.Range('E44').Value = '  -7.97%  '
$ws.Range('D45').Value = '10.17'
$ws.Range('E45').Value = '  -1.23%  '
$ws.Range('D46').Value = '17.34'
$ws.Range('E46').Value = '  -2.64%  '
$ws.Range('E47').Value = '  -3.59%  '
$ws.Range('D48').Value = '54.36'
$ws.Range('E48').Value = '  -3.53%  '
$ws.Range('D49').Value = '2.530.05'
$ws.Range('E49').Value = '  -0.47%  '
$ws.Range('E50').Value = '  -2.40%  '
$ws.Range('B51').Value = 'BitcoinSV'
$ws.Range('C51').Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range('D51').Value = '72.74'
$ws.Range('E51').Value = '  +2.46%  '

$editRange.ClearFormats()
